# Apply the tracked changes from the commit:
#  - Update the two description strings in column B
#  - Widen column B
#  - Reset rows 2 and 3 back to the sheet's default (non-custom) height
#  - Move the active selection to B3
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update cell text (shared-string content changes)
$ws.Range("B2").Value = "Colour of ball drawn from Bag"
$ws.Range("B3").Value = "Bag number from which ball is drawn"

# 2) Widen column B (ColumnWidth is in characters; this lands on the
#    character-width bucket closest to the target 30.92 OOXML width)
$ws.Columns("B").ColumnWidth = 30.14

# 3) Rows 2 and 3 go back to the implicit/default row height (12.8pt,
#    same as sheetFormatPr's defaultRowHeight) instead of the custom 14.65pt
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()

# 4) Move the selection/active cell to B3
$ws.Range("B3").Select() | Out-Null
